$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_ETR"

$values = @(
    0,
    22.3629700097761,
    24.1070977083343,
    2.70143560365249,
    0,
    27.1036868139497,
    25.7633956084017,
    7.4336205941199,
    0,
    22.8897305934162,
    30.1022447923085,
    1.55989780120439,
    0,
    32.1993956629253,
    22.7342285164783,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
